# Update the Jengibre price-history sheet: refresh weekly data
# (dates + volume/price columns), appending the two new rows
# pulled in by this week's run and extending the sheet to row 25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows (2-23): refresh Fecha / Volumen / Precio columns ---
$ws.Cells.Item(2, 4).Value = 44326
$ws.Cells.Item(2, 10).Value = 460
$ws.Cells.Item(2, 11).Value = 25000
$ws.Cells.Item(2, 12).Value = 26000
$ws.Cells.Item(2, 13).Value = 25500
$ws.Cells.Item(2, 16).Value = 1962

$ws.Cells.Item(3, 4).Value = 44383
$ws.Cells.Item(3, 10).Value = 200
$ws.Cells.Item(3, 11).Value = 17000
$ws.Cells.Item(3, 12).Value = 18000
$ws.Cells.Item(3, 13).Value = 17500
$ws.Cells.Item(3, 16).Value = 1346

$ws.Cells.Item(4, 4).Value = 44445
$ws.Cells.Item(4, 10).Value = 600
$ws.Cells.Item(4, 11).Value = 13000
$ws.Cells.Item(4, 12).Value = 14000
$ws.Cells.Item(4, 13).Value = 13500
$ws.Cells.Item(4, 16).Value = 1038

$ws.Cells.Item(5, 4).Value = 44333
$ws.Cells.Item(5, 10).Value = 440
$ws.Cells.Item(5, 11).Value = 24000
$ws.Cells.Item(5, 12).Value = 25000
$ws.Cells.Item(5, 13).Value = 24500
$ws.Cells.Item(5, 16).Value = 1885

$ws.Cells.Item(6, 4).Value = 44323
$ws.Cells.Item(6, 10).Value = 460
$ws.Cells.Item(6, 11).Value = 25000
$ws.Cells.Item(6, 12).Value = 26000
$ws.Cells.Item(6, 13).Value = 25500
$ws.Cells.Item(6, 16).Value = 1962

$ws.Cells.Item(7, 4).Value = 44242
$ws.Cells.Item(7, 10).Value = 400
$ws.Cells.Item(7, 11).Value = 44000
$ws.Cells.Item(7, 12).Value = 45000
$ws.Cells.Item(7, 13).Value = 44500
$ws.Cells.Item(7, 16).Value = 3423

$ws.Cells.Item(8, 4).Value = 44260
$ws.Cells.Item(8, 10).Value = 400
$ws.Cells.Item(8, 11).Value = 37000
$ws.Cells.Item(8, 12).Value = 38000
$ws.Cells.Item(8, 13).Value = 37500
$ws.Cells.Item(8, 16).Value = 2885

$ws.Cells.Item(9, 4).Value = 44418
$ws.Cells.Item(9, 10).Value = 500
$ws.Cells.Item(9, 11).Value = 14000
$ws.Cells.Item(9, 12).Value = 15000
$ws.Cells.Item(9, 13).Value = 14500
$ws.Cells.Item(9, 16).Value = 1115

$ws.Cells.Item(10, 4).Value = 44312
$ws.Cells.Item(10, 10).Value = 400
$ws.Cells.Item(10, 11).Value = 26000
$ws.Cells.Item(10, 12).Value = 27000
$ws.Cells.Item(10, 13).Value = 26500
$ws.Cells.Item(10, 16).Value = 2038

$ws.Cells.Item(11, 4).Value = 44410
$ws.Cells.Item(11, 10).Value = 600
$ws.Cells.Item(11, 11).Value = 14000
$ws.Cells.Item(11, 12).Value = 15000
$ws.Cells.Item(11, 13).Value = 14500
$ws.Cells.Item(11, 16).Value = 1115

$ws.Cells.Item(12, 4).Value = 44344
$ws.Cells.Item(12, 10).Value = 400
$ws.Cells.Item(12, 11).Value = 18500
$ws.Cells.Item(12, 12).Value = 19000
$ws.Cells.Item(12, 13).Value = 18750
$ws.Cells.Item(12, 16).Value = 1442

$ws.Cells.Item(13, 4).Value = 44419
$ws.Cells.Item(13, 10).Value = 600
$ws.Cells.Item(13, 11).Value = 14000
$ws.Cells.Item(13, 12).Value = 15000
$ws.Cells.Item(13, 13).Value = 14500
$ws.Cells.Item(13, 16).Value = 1115

$ws.Cells.Item(14, 4).Value = 44249
$ws.Cells.Item(14, 10).Value = 400
$ws.Cells.Item(14, 11).Value = 42000
$ws.Cells.Item(14, 12).Value = 43000
$ws.Cells.Item(14, 13).Value = 42500
$ws.Cells.Item(14, 16).Value = 3269

$ws.Cells.Item(15, 4).Value = 44309
$ws.Cells.Item(15, 10).Value = 400
$ws.Cells.Item(15, 11).Value = 26000
$ws.Cells.Item(15, 12).Value = 27000
$ws.Cells.Item(15, 13).Value = 26500
$ws.Cells.Item(15, 16).Value = 2038

$ws.Cells.Item(16, 4).Value = 44365
$ws.Cells.Item(16, 10).Value = 500
$ws.Cells.Item(16, 11).Value = 19500
$ws.Cells.Item(16, 12).Value = 20000
$ws.Cells.Item(16, 13).Value = 19750
$ws.Cells.Item(16, 16).Value = 1519

$ws.Cells.Item(17, 4).Value = 44335
$ws.Cells.Item(17, 10).Value = 480
$ws.Cells.Item(17, 11).Value = 24500
$ws.Cells.Item(17, 12).Value = 25000
$ws.Cells.Item(17, 13).Value = 24750
$ws.Cells.Item(17, 16).Value = 1904

$ws.Cells.Item(18, 4).Value = 44426
$ws.Cells.Item(18, 10).Value = 460
$ws.Cells.Item(18, 11).Value = 14000
$ws.Cells.Item(18, 12).Value = 15000
$ws.Cells.Item(18, 13).Value = 14500
$ws.Cells.Item(18, 16).Value = 1115

$ws.Cells.Item(19, 4).Value = 44414
$ws.Cells.Item(19, 10).Value = 500
$ws.Cells.Item(19, 11).Value = 14000
$ws.Cells.Item(19, 12).Value = 15000
$ws.Cells.Item(19, 13).Value = 14500
$ws.Cells.Item(19, 16).Value = 1115

$ws.Cells.Item(20, 4).Value = 44379
$ws.Cells.Item(20, 10).Value = 600
$ws.Cells.Item(20, 11).Value = 17000
$ws.Cells.Item(20, 12).Value = 18000
$ws.Cells.Item(20, 13).Value = 17500
$ws.Cells.Item(20, 16).Value = 1346

$ws.Cells.Item(21, 4).Value = 44428
$ws.Cells.Item(21, 10).Value = 480
$ws.Cells.Item(21, 11).Value = 14000
$ws.Cells.Item(21, 12).Value = 15000
$ws.Cells.Item(21, 13).Value = 14500
$ws.Cells.Item(21, 16).Value = 1115

$ws.Cells.Item(22, 4).Value = 44442
$ws.Cells.Item(22, 10).Value = 460
$ws.Cells.Item(22, 11).Value = 14000
$ws.Cells.Item(22, 12).Value = 15000
$ws.Cells.Item(22, 13).Value = 14500
$ws.Cells.Item(22, 16).Value = 1115

$ws.Cells.Item(23, 4).Value = 44435
$ws.Cells.Item(23, 10).Value = 480
$ws.Cells.Item(23, 11).Value = 13000
$ws.Cells.Item(23, 12).Value = 14000
$ws.Cells.Item(23, 13).Value = 13500
$ws.Cells.Item(23, 16).Value = 1038

# --- New rows (24-25): same market/category, new dates & prices ---
$ws.Cells.Item(24, 1).Value = 8
$ws.Cells.Item(24, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(24, 3).Value = "Coquimbo"
$ws.Cells.Item(24, 4).Value = 44412
$ws.Cells.Item(24, 5).Value = 4
$ws.Cells.Item(24, 6).Value = 100114007
$ws.Cells.Item(24, 7).Value = "Jengibre"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 600
$ws.Cells.Item(24, 11).Value = 14000
$ws.Cells.Item(24, 12).Value = 15000
$ws.Cells.Item(24, 13).Value = 14500
$ws.Cells.Item(24, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(24, 15).Value = "Perú"
$ws.Cells.Item(24, 16).Value = 1115
$ws.Cells.Item(24, 17).Value = 13
$ws.Cells.Item(24, 18).Value = "Hortaliza"

$ws.Cells.Item(25, 1).Value = 8
$ws.Cells.Item(25, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(25, 3).Value = "Coquimbo"
$ws.Cells.Item(25, 4).Value = 44400
$ws.Cells.Item(25, 5).Value = 4
$ws.Cells.Item(25, 6).Value = 100114007
$ws.Cells.Item(25, 7).Value = "Jengibre"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 600
$ws.Cells.Item(25, 11).Value = 15000
$ws.Cells.Item(25, 12).Value = 16000
$ws.Cells.Item(25, 13).Value = 15500
$ws.Cells.Item(25, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(25, 15).Value = "Perú"
$ws.Cells.Item(25, 16).Value = 1192
$ws.Cells.Item(25, 17).Value = 13
$ws.Cells.Item(25, 18).Value = "Hortaliza"

# --- Match formatting of the new rows to the existing data rows ---
$ws.Range("A23:R23").Copy()
$ws.Range("A24:R25").PasteSpecial(-4122)
$excel.CutCopyMode = 0
